# Sprint 2 burndown update (stand-up on Nov 4 / commit notes "Nov 30"):
#  - Tasks Left actuals for the last 3 tracked days drop from 20 to 7
#  - Chart title date is advanced
#  - Selection left on D9 as the sheet was last viewed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Tasks Left" (column C) actuals for Oct 30 / Nov 1 / Nov 4 flatten at 7
$ws.Range("C4").Value = 7
$ws.Range("C5").Value = 7
$ws.Range("C6").Value = 7

# Chart title: "...as of October 28" -> "...as of November 4"
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.ChartTitle.Text = "Burndown Chart  for Sprint 2 as of November 4"

# Active cell/selection on the sheet moved from C9 to D9
$ws.Range("D9").Select()
